$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in team record values for each data row (2-43)
for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 30).Value = 85   # AD = col 30
    $ws.Cells.Item($row, 31).Value = 77   # AE = col 31
    $ws.Cells.Item($row, 32).Value = 0    # AF = col 32
}
